# WordCount.xlsx: Data Warehousing complete, now all of Data Analytics completed
# -> fill in the actual word counts for "Data Mining" (row 6) and
#    "Data Warehousing" (row 7); the Total (C11, =SUM(C2:C10)) recalculates
#    automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = 397
$ws.Range("C7").Value = 282

# Column contents got wider (e.g. "397"/"282" vs "0"), so the best-fit
# columns re-size themselves, same as Excel auto-adjusting on edit.
$ws.Columns.Item(1).ColumnWidth = 27.3
$ws.Columns.Item(2).ColumnWidth = 20
$ws.Columns.Item(3).ColumnWidth = 17.1

# Cursor ends up on the next cell of the entry column after the edits.
$ws.Range("C8").Select()
